$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (values like "96.32"). Each cell is set
# individually since this runtime only honors NumberFormat on the first area
# of a multi-area Range.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.600.02"
$ws.Range("E2").Value = "  -1.43%  "

# Row 3
$ws.Range("D3").Value = "2.289.45"
$ws.Range("E3").Value = "  +0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "96.32"
$ws.Range("E5").Value = "  -1.59%  "

# Row 6
$ws.Range("D6").Value = "268.09"
$ws.Range("E6").Value = "  -2.62%  "

# Row 7
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -2.09%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  -5.55%  "

# Row 10
$ws.Range("D10").Value = "45.19"
$ws.Range("E10").Value = "  -5.92%  "

# Row 11
$ws.Range("E11").Value = "  -0.76%  "

# Row 12
$ws.Range("D12").Value = "7.82"
$ws.Range("E12").Value = "  -4.72%  "

# Row 13
$ws.Range("E13").Value = "  -0.13%  "

# Row 14
$ws.Range("D14").Value = "2.632.51"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15
$ws.Range("E15").Value = "  -3.39%  "

# Row 16
$ws.Range("D16").Value = "0.851"
$ws.Range("E16").Value = "  +0.88%  "

# Row 17
$ws.Range("D17").Value = "2.293.97"
$ws.Range("E17").Value = "  +1.63%  "

# Row 18
$ws.Range("D18").Value = "43.563.27"
$ws.Range("E18").Value = "  -1.46%  "

# Row 19
$ws.Range("E19").Value = "  +1.19%  "

# Row 20
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  -0.97%  "

# Row 21
$ws.Range("D21").Value = "71.95"
$ws.Range("E21").Value = "  +1.17%  "

# Row 22
$ws.Range("E22").Value = "  +4.38%  "

# Row 23
$ws.Range("D23").Value = "233.00"
$ws.Range("E23").Value = "  -1.02%  "

# Row 24
$ws.Range("D24").Value = "9.17"
$ws.Range("E24").Value = "  -10.62%  "

# Row 25
$ws.Range("E25").Value = "  -0.11%  "

# Row 26
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").Value = "  -1.86%  "

# Row 27
$ws.Range("D27").Value = "11.20"
$ws.Range("E27").Value = "  -3.31%  "

# Row 28
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").Value = "3.46"
$ws.Range("E28").Value = "  +3.32%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "40.48"
$ws.Range("E29").Value = "  +2.03%  "

# Row 30
$ws.Range("E30").Value = "  +0.31%  "

# Row 31
$ws.Range("D31").Value = "175.34"
$ws.Range("E31").Value = "  +0.93%  "

# Row 32
$ws.Range("D32").Value = "21.91"
$ws.Range("E32").Value = "  +2.43%  "

# Row 33
$ws.Range("D33").Value = "0.0883"
$ws.Range("E33").Value = "  -4.53%  "

# Row 34
$ws.Range("E34").Value = "  -5.79%  "

# Row 35
$ws.Range("E35").Value = "  +0.39%  "

# Row 36
$ws.Range("E36").Value = "  -5.12%  "

# Row 37
$ws.Range("E37").Value = "  -0.95%  "

# Row 38
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  -0.96%  "

# Row 39
$ws.Range("E39").Value = "  -6.35%  "

# Row 40
$ws.Range("D40").Value = "0.236"
$ws.Range("E40").Value = "  -8.30%  "

# Row 41
$ws.Range("E41").Value = "  +5.49%  "

# Row 42
$ws.Range("E42").Value = "  -3.02%  "

# Row 43
$ws.Range("D43").Value = "1.35"
$ws.Range("E43").Value = "  +14.57%  "

# Row 44
$ws.Range("D44").Value = "63.95"
$ws.Range("E44").Value = "  +1.58%  "

# Row 45
$ws.Range("E45").Value = "  +2.38%  "

# Row 46
$ws.Range("E46").Value = "  -5.18%  "

# Row 47
$ws.Range("E47").Value = "  -0.79%  "

# Row 48
$ws.Range("D48").Value = "98.23"
$ws.Range("E48").Value = "  -2.37%  "

# Row 49
$ws.Range("E49").Value = "  -0.22%  "

# Row 50
$ws.Range("D50").Value = "2.512.47"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51
$ws.Range("E51").Value = "  -1.98%  "
